$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3496

$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -4716

$ws.Range("H129").Value = 1200.6129
$ws.Range("I129").Value = 900.5
$ws.Range("J129").Value = 1210.6167
$ws.Range("K129").Value = 2701.5
$ws.Range("L129").Value = 3631.8501
$ws.Range("M129").Value = 2298.5
$ws.Range("N129").Value = -13631.8501

$ws.Range("H135").Value = 16130218
$ws.Range("I135").Value = 20834322
$ws.Range("J135").Value = 1862.2858
$ws.Range("K135").Value = 187508898
$ws.Range("L135").Value = 16760.5722
$ws.Range("M135").Value = -187506363
$ws.Range("N135").Value = -21830.5722

$ws.Range("H137").Value = 1939644.9
$ws.Range("I137").Value = 3473418.8
$ws.Range("K137").Value = 10420256.4
$ws.Range("M137").Value = -10417706.4

$ws.Range("H138").Value = 4000.8162
$ws.Range("I138").Value = 3259.75
$ws.Range("J138").Value = 4283.127
$ws.Range("K138").Value = 9779.25
$ws.Range("L138").Value = 12849.381
$ws.Range("M138").Value = -4639.25
$ws.Range("N138").Value = -23129.381

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3952.34
$ws.Range("I32").Value = 3952.34
$ws.Range("K32").Value = 3952.34
$ws.Range("M32").Value = -3665.34

$ws.Range("H61").Value = 9011353
$ws.Range("I61").Value = 13334658
$ws.Range("K61").Value = 13334658
$ws.Range("M61").Value = -13334446

$ws.Range("H74").Value = 12822338
$ws.Range("I74").Value = 1061.3334
$ws.Range("K74").Value = 1061.3334
$ws.Range("M74").Value = -187.3334

$ws.Range("H77").Value = 12822338
$ws.Range("I77").Value = 1061.3334
$ws.Range("K77").Value = 5306.666999999999
$ws.Range("M77").Value = -938.6669999999995

$ws.Range("H80").Value = 19186.5
$ws.Range("J80").Value = 19186.5
$ws.Range("L80").Value = 19186.5
$ws.Range("N80").Value = -21182.5

$ws.Range("H83").Value = 19186.5
$ws.Range("J83").Value = 19186.5
$ws.Range("L83").Value = 57559.5
$ws.Range("N83").Value = -67543.5

$ws.Range("H102").Value = 4261.1113
$ws.Range("I102").Value = 3608.3333
$ws.Range("J102").Value = 5566.6665
$ws.Range("K102").Value = 3608.3333
$ws.Range("L102").Value = 5566.6665
$ws.Range("M102").Value = -1986.3333
$ws.Range("N102").Value = -8810.666499999999

$ws.Range("H110").Value = 57006.875
$ws.Range("I110").Value = 82583.82000000001
$ws.Range("J110").Value = 737.6
$ws.Range("K110").Value = 82583.82000000001
$ws.Range("L110").Value = 737.6
$ws.Range("M110").Value = -80538.82000000001
$ws.Range("N110").Value = -4827.6

$ws.Range("H132").Value = 1071088.4
$ws.Range("I132").Value = 2005.7548
$ws.Range("J132").Value = 4053266.5
$ws.Range("K132").Value = 6017.2644
$ws.Range("L132").Value = 12159799.5
$ws.Range("M132").Value = -3487.2644
$ws.Range("N132").Value = -12164859.5

$ws.Range("H136").Value = 9011353
$ws.Range("I136").Value = 13334658
$ws.Range("K136").Value = 40003974
$ws.Range("M136").Value = -40001424

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11365617
$ws.Range("I105").Value = 13890532
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 13890532
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -13888785
$ws.Range("N105").Value = -6994

$ws.Range("H107").Value = 64152.812
$ws.Range("I107").Value = 72896.07000000001
$ws.Range("J107").Value = 2950
$ws.Range("K107").Value = 72896.07000000001
$ws.Range("L107").Value = 2950
$ws.Range("M107").Value = -70976.07000000001
$ws.Range("N107").Value = -6790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 18750.25
$ws.Range("J4").Value = 18750.25
$ws.Range("L4").Value = 18750.25
$ws.Range("N4").Value = -18974.25

$ws.Range("H31").Value = 5692.5
$ws.Range("I31").Value = 2100.2354
$ws.Range("J31").Value = 9431.388000000001
$ws.Range("K31").Value = 2100.2354
$ws.Range("L31").Value = 9431.388000000001
$ws.Range("M31").Value = -1805.2354
$ws.Range("N31").Value = -10021.388

$ws.Range("H34").Value = 5692.5
$ws.Range("I34").Value = 2100.2354
$ws.Range("J34").Value = 9431.388000000001
$ws.Range("K34").Value = 2100.2354
$ws.Range("L34").Value = 9431.388000000001
$ws.Range("M34").Value = -1898.2354
$ws.Range("N34").Value = -9835.388000000001

$ws.Range("H58").Value = 1880.52
$ws.Range("I58").Value = 1736.7222
$ws.Range("J58").Value = 2250.2856
$ws.Range("K58").Value = 1736.7222
$ws.Range("L58").Value = 2250.2856
$ws.Range("M58").Value = -1533.7222
$ws.Range("N58").Value = -2656.2856

$ws.Range("H136").Value = 1880.52
$ws.Range("I136").Value = 1736.7222
$ws.Range("J136").Value = 2250.2856
$ws.Range("K136").Value = 5210.1666
$ws.Range("L136").Value = 6750.8568
$ws.Range("M136").Value = -2660.1666
$ws.Range("N136").Value = -11850.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1498.75
$ws.Range("I60").Value = 597
$ws.Range("J60").Value = 2142.8572
$ws.Range("K60").Value = 1791
$ws.Range("L60").Value = 6428.571599999999
$ws.Range("M60").Value = -1540
$ws.Range("N60").Value = -6930.571599999999

$ws.Range("H81").Value = 6363.273
$ws.Range("J81").Value = 14749.25
$ws.Range("L81").Value = 44247.75
$ws.Range("N81").Value = -46493.75

$ws.Range("H84").Value = 6363.273
$ws.Range("J84").Value = 14749.25
$ws.Range("L84").Value = 132743.25
$ws.Range("N84").Value = -143975.25

$ws.Range("H107").Value = 2241
$ws.Range("I107").Value = 410
$ws.Range("J107").Value = 2764.1428
$ws.Range("K107").Value = 1230
$ws.Range("L107").Value = 8292.428400000001
$ws.Range("M107").Value = 690
$ws.Range("N107").Value = -12132.4284

$ws.Range("H113").Value = 634.95654
$ws.Range("I113").Value = 616.92
$ws.Range("J113").Value = 656.4286
$ws.Range("K113").Value = 1850.76
$ws.Range("L113").Value = 1969.2858
$ws.Range("M113").Value = 319.2400000000002
$ws.Range("N113").Value = -6309.2858

$ws.Range("H137").Value = 36049.03
$ws.Range("I137").Value = 6284.0835
$ws.Range("J137").Value = 115422.22
$ws.Range("K137").Value = 18852.2505
$ws.Range("L137").Value = 346266.66
$ws.Range("M137").Value = -13752.2505
$ws.Range("N137").Value = -356466.66

$ws.Range("H140").Value = 1443.8363
$ws.Range("I140").Value = 974.2564
$ws.Range("J140").Value = 2588.4375
$ws.Range("K140").Value = 2922.7692
$ws.Range("L140").Value = 7765.3125
$ws.Range("M140").Value = 2257.2308
$ws.Range("N140").Value = -18125.3125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 75.8
$ws.Range("I2").Value = 75.8
$ws.Range("K2").Value = 75.8
$ws.Range("M2").Value = 37.2

$ws.Range("H43").Value = 6340.846
$ws.Range("I43").Value = 679.75
$ws.Range("J43").Value = 15398.6
$ws.Range("K43").Value = 679.75
$ws.Range("L43").Value = 15398.6
$ws.Range("M43").Value = -528.75
$ws.Range("N43").Value = -15700.6

$ws.Range("H80").Value = 35215228
$ws.Range("I80").Value = 87833336
$ws.Range("J80").Value = 136488.67
$ws.Range("K80").Value = 87833336
$ws.Range("L80").Value = 136488.67
$ws.Range("M80").Value = -87832338
$ws.Range("N80").Value = -138484.67

$ws.Range("H83").Value = 35215228
$ws.Range("I83").Value = 87833336
$ws.Range("J83").Value = 136488.67
$ws.Range("K83").Value = 439166680
$ws.Range("L83").Value = 682443.3500000001
$ws.Range("M83").Value = -439161688
$ws.Range("N83").Value = -692427.3500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9911.777
$ws.Range("J2").Value = 9911.777
$ws.Range("L2").Value = 9911.777
$ws.Range("N2").Value = -10135.777

$ws.Range("H61").Value = 6101
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 6101
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 6101
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -6505

$ws.Range("H113").Value = 6101
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 6101
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6101
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10441

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1672666.6
$ws.Range("I14").Value = 8000
$ws.Range("J14").Value = 2505000
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 2505000
$ws.Range("M14").Value = -7832
$ws.Range("N14").Value = -2505336

$ws.Range("H98").Value = 95000
$ws.Range("J98").Value = 95000
$ws.Range("L98").Value = 95000
$ws.Range("N98").Value = -100990

$ws.Range("H100").Value = 796.1905
$ws.Range("I100").Value = 629.7273
$ws.Range("K100").Value = 1259.4546
$ws.Range("M100").Value = -718.4546

$ws.Range("H105").Value = 52871.668
$ws.Range("J105").Value = 52871.668
$ws.Range("L105").Value = 52871.668
$ws.Range("N105").Value = -59859.668

$ws.Range("H108").Value = 70313
$ws.Range("J108").Value = 70313
$ws.Range("L108").Value = 70313
$ws.Range("N108").Value = -77993
